$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Current (before) state:
# Row 6: A6=107258607, Q6=404755.5111078721, R6=7063764.822795196, M6=empty
# Row 7: A7=107258628, Q7=404588.0690095468, R7=7064520.029476826, M7=empty
# Row 8: A8=107258608, Q8=404465.3900776547, R8=7064504.653031247, M8="färska spår"
#
# Target (after) state is a cyclic rotation: new row6 = old row7, new row7 = old row8, new row8 = old row6
# for columns A, M, Q, R.

$ws.Range("A6").Value = 107258628
$ws.Range("Q6").Value = 404588.0690095468
$ws.Range("R6").Value = 7064520.029476826

$ws.Range("A7").Value = 107258608
$ws.Range("M7").Value = "färska spår"
$ws.Range("Q7").Value = 404465.3900776547
$ws.Range("R7").Value = 7064504.653031247

$ws.Range("A8").Value = 107258607
$ws.Range("M8").Value = ""
$ws.Range("Q8").Value = 404755.5111078721
$ws.Range("R8").Value = 7063764.822795196
